$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 value: "moneter" -> "ekonomi,moneter"
$ws.Range("B1").Value = "ekonomi,moneter"

# Update B8 value: "kepemimpinan,sekolah,budaya mutu" -> "kepemimpinan,sekolah,budayamutu"
$ws.Range("B8").Value = "kepemimpinan,sekolah,budayamutu"

# Remove column C entirely (data + formatting) - tag-along "Keterangan" column no longer needed
$ws.Columns.Item(3).Delete()

# Update selection to B10
$ws.Range("B10").Select()
